# Updates the cryptocurrency price/volume table on Sheet1 to reflect
# the latest scrape (GitHub Actions scheduled refresh).
# Rows 39/40 and 42/43 also swap rank order (Stacks/PEPE, Fetch.AI/CoreDAO).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text, even when it looks like a number
# (e.g. '1.00', '556.16'), so the cell keeps its original text formatting
# instead of being auto-converted into a numeric value by Excel.
function Set-TextValue($address, $text) {
    $cell = $ws.Range($address)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = "65.784.73"
$ws.Range("E2").Value = "  -5.58%  "

# Row 3
$ws.Range("D3").Value = "3.273.35"
$ws.Range("E3").Value = "  -6.38%  "

# Row 4
$ws.Range("E4").Value = "  -0.14%  "

# Row 5
Set-TextValue "D5" "556.16"
$ws.Range("E5").Value = "  -3.77%  "

# Row 6
Set-TextValue "D6" "183.68"
$ws.Range("E6").Value = "  -4.48%  "

# Row 7
Set-TextValue "D7" "1.00"
$ws.Range("E7").Value = "  +0.02%  "

# Row 8
Set-TextValue "D8" "0.588"
$ws.Range("E8").Value = "  -4.03%  "

# Row 9
$ws.Range("D9").Value = "3.263.11"
$ws.Range("E9").Value = "  -6.37%  "

# Row 10
Set-TextValue "D10" "0.185"
$ws.Range("E10").Value = "  -10.04%  "

# Row 11
Set-TextValue "D11" "0.583"
$ws.Range("E11").Value = "  -6.13%  "

# Row 12
Set-TextValue "D12" "47.16"
$ws.Range("E12").Value = "  -8.37%  "

# Row 13
Set-TextValue "D13" "0.0000266"
$ws.Range("E13").Value = "  -7.20%  "

# Row 14
Set-TextValue "D14" "644.45"
$ws.Range("E14").Value = "  -0.17%  "

# Row 15
Set-TextValue "D15" "8.65"
$ws.Range("E15").Value = "  -5.49%  "

# Row 16
$ws.Range("D16").Value = "3.796.07"
$ws.Range("E16").Value = "  -6.42%  "

# Row 17
Set-TextValue "D17" "18.03"
$ws.Range("E17").Value = "  -1.69%  "

# Row 18
$ws.Range("D18").Value = "65.701.87"
$ws.Range("E18").Value = "  -5.79%  "

# Row 19
$ws.Range("E19").Value = "  -3.31%  "

# Row 20
$ws.Range("D20").Value = "3.262.87"
$ws.Range("E20").Value = "  -6.54%  "

# Row 21
Set-TextValue "D21" "11.33"
$ws.Range("E21").Value = "  -8.77%  "

# Row 22
Set-TextValue "D22" "0.902"
$ws.Range("E22").Value = "  -4.98%  "

# Row 23
Set-TextValue "D23" "18.30"
$ws.Range("E23").Value = "  +1.46%  "

# Row 24
Set-TextValue "D24" "108.12"
$ws.Range("E24").Value = "  +9.32%  "

# Row 25
Set-TextValue "D25" "4.92"
$ws.Range("E25").Value = "  -7.89%  "

# Row 26
Set-TextValue "D26" "3.97"
$ws.Range("E26").Value = "  -7.44%  "

# Row 27
Set-TextValue "D27" "2.67"
$ws.Range("E27").Value = "  -7.18%  "

# Row 28
Set-TextValue "D28" "9.60"
$ws.Range("E28").Value = "  -4.68%  "

# Row 29
Set-TextValue "D29" "8.64"
$ws.Range("E29").Value = "  -8.02%  "

# Row 30
Set-TextValue "D30" "30.21"
$ws.Range("E30").Value = "  -7.55%  "

# Row 31
Set-TextValue "D31" "3.95"
$ws.Range("E31").Value = "  -6.66%  "

# Row 32
Set-TextValue "D32" "6.27"
$ws.Range("E32").Value = "  -6.82%  "

# Row 33
Set-TextValue "D33" "11.04"
$ws.Range("E33").Value = "  -5.37%  "

# Row 34
Set-TextValue "D34" "0.105"
$ws.Range("E34").Value = "  -4.67%  "

# Row 35
$ws.Range("D35").Value = "3.753.04"
$ws.Range("E35").Value = "  +1.34%  "

# Row 36
Set-TextValue "D36" "57.47"
$ws.Range("E36").Value = "  -6.54%  "

# Row 37
Set-TextValue "D37" "1.00"
$ws.Range("E37").Value = "  -0.08%  "

# Row 38
Set-TextValue "D38" "517.95"
$ws.Range("E38").Value = "  -7.82%  "

# Row 39
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D39" "3.40"
$ws.Range("E39").Value = "  -6.45%  "

# Row 40
$ws.Range("B40").Value = "PEPE"
$ws.Range("C40").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D40").Value = "0.0₃0735"
$ws.Range("E40").Value = "  -6.89%  "

# Row 41
$ws.Range("E41").Value = "  -1.95%  "

# Row 42
$ws.Range("B42").Value = "Fetch.AI"
$ws.Range("C42").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D42" "2.73"
$ws.Range("E42").Value = "  -5.88%  "

# Row 43
$ws.Range("B43").Value = "CoreDAO"
$ws.Range("C43").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
Set-TextValue "D43" "3.41"
$ws.Range("E43").Value = "  -13.75%  "

# Row 44
Set-TextValue "D44" "32.88"
$ws.Range("E44").Value = "  -4.04%  "

# Row 45
Set-TextValue "D45" "0.336"
$ws.Range("E45").Value = "  -10.19%  "

# Row 46
Set-TextValue "D46" "0.0413"
$ws.Range("E46").Value = "  -6.88%  "

# Row 47
Set-TextValue "D47" "3.21"
$ws.Range("E47").Value = "  -4.89%  "

# Row 48
$ws.Range("E48").Value = "  -4.36%  "

# Row 49
Set-TextValue "D49" "2.61"
$ws.Range("E49").Value = "  -8.38%  "

# Row 50
$ws.Range("E50").Value = "  -0.05%  "

# Row 51
$ws.Range("E51").Value = "  +1.40%  "
